# Re-orders the 20 data rows (rows 2-21) of the active worksheet according
# to the permutation observed between the "before" and "after" versions of
# the workbook. Row 1 (header) is left untouched. For each target row, the
# full row (columns A:AY) is replaced with the original content that used
# to live in the corresponding source row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row number -> source row number (1-based row numbers on the sheet)
$mapping = @{
  2  = 10
  3  = 11
  4  = 2
  5  = 12
  6  = 3
  7  = 4
  8  = 5
  9  = 13
  10 = 14
  11 = 15
  12 = 6
  13 = 16
  14 = 17
  15 = 18
  16 = 19
  17 = 20
  18 = 21
  19 = 7
  20 = 8
  21 = 9
}

$firstRow = 2
$lastRow = 21
$lastCol = "AY"

# Snapshot every data row's full contents (columns A:AY) before making any
# changes, since the re-ordering is a single large cycle and writing in
# place would otherwise clobber data that is still needed.
$snapshots = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rng = $ws.Range("A${r}:${lastCol}${r}")
    $snapshots[$r] = $rng.Value2
}

# Columns Y and AA hold plain date text (e.g. "2021-08-24"). Force them to
# a text number format on every data row so that re-assigning the snapshot
# values does not get auto-converted into Excel date serial numbers.
$ws.Range("Y${firstRow}:Y${lastRow}").NumberFormat = "@"
$ws.Range("AA${firstRow}:AA${lastRow}").NumberFormat = "@"

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    $data = $snapshots[$source]
    $dst = $ws.Range("A${target}:${lastCol}${target}")
    $dst.Value2 = $data
}
